$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update latest-period value for "Employment share by occupation" row ---
$ws.Range("C3").Value = "Jul 2021 - Jun 2022 (11/10/22)"

# --- Append new data-source rows ---
$ws.Range("A7").Value = "Enterprise by employment size"
$ws.Range("B7").Value = "<a href='https://www.nomisweb.co.uk/datasets/idbrent'>ONS UK Business Count</a>"
$ws.Range("C7").Value = "Oct 2021 - Sept 2022 (28/09/22)"
$ws.Range("D7").Value = "Oct 2022 - Sept 2023 (03/10/23)"

$ws.Range("A8").Value = "Key Stage 4 (KS4) destinations "
$ws.Range("B8").Value = "<a href = 'https://explore-education-statistics.service.gov.uk/data-tables/permalink/c9f44a09-4239-40d6-8f07-87c2b97fc5fc'>National Pupil Database</a>"
$ws.Range("C8").Value = "Aug 2019 -  Jul 2020 (2019 leavers) (21/10/21)"
$ws.Range("D8").Value = "Aug 2020 - Jul 2021 (2020 leavers) (20/10/22)"

$ws.Range("A9").Value = "Key Stage 5 (KS5) destinations "
$ws.Range("B9").Value = "<a href = 'https://explore-education-statistics.service.gov.uk/data-tables/permalink/6ed2058c-1ff4-4e13-b167-3b15bb6a0675'>National Pupil Database</a>"
$ws.Range("C9").Value = "Aug 2019 - Jul 2020 (2019 leavers) (09/12/21)"
$ws.Range("D9").Value = "Aug 2020 - Jul 2021 (2020 leavers) (20/10/22)"

# --- Widen column D to fit the newer, longer "next period" strings ---
$ws.Columns.Item(4).ColumnWidth = 37.9

# --- Restore the selection to where the author left off editing ---
[void]$ws.Range("C14").Select()
